# Weekly fruit/vegetable price update for "Pepino dulce" (La Palmera de La Serena).
# A new week's price block (44628 = 2022-03-08) is written into rows 304-306
# (Primera/Segunda/Tercera only - no "Especial" grade reported that week),
# and the previous week's block (44335 = 2021-05-19, Especial/Primera/Segunda/
# Tercera) that used to occupy rows 304-307 is preserved, shifted down into
# rows 307-310.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = "YYYY-MM-DD HH:MM:SS"

# Common, unchanged columns shared by every row in this block.
$mercadoId = 8
$mercado   = "Terminal La Palmera de La Serena"
$region    = "Coquimbo"
$codreg    = 4
$catId     = 100112043
$categoria = "Pepino dulce"
$variedad  = "Cultivar IV Región"
$unidad    = "`$/bandeja 18 kilos"
$origen    = "Provincia de Limarí"
$clasif    = "Hortaliza"

function Set-PrecioRow($Row, $Fecha, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg) {
    $ws.Cells.Item($Row, 1).Value = $mercadoId
    $ws.Cells.Item($Row, 2).Value = $mercado
    $ws.Cells.Item($Row, 3).Value = $region

    $ws.Cells.Item($Row, 4).NumberFormat = $dateFormat
    $ws.Cells.Item($Row, 4).Value = $Fecha

    $ws.Cells.Item($Row, 5).Value = $codreg
    $ws.Cells.Item($Row, 6).Value = $catId
    $ws.Cells.Item($Row, 7).Value = $categoria
    $ws.Cells.Item($Row, 8).Value = $variedad
    $ws.Cells.Item($Row, 9).Value = $Calidad
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMin
    $ws.Cells.Item($Row, 12).Value = $PrecioMax
    $ws.Cells.Item($Row, 13).Value = $PrecioProm
    $ws.Cells.Item($Row, 14).Value = $unidad
    $ws.Cells.Item($Row, 15).Value = $origen
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
    $ws.Cells.Item($Row, 17).Value = 18
    $ws.Cells.Item($Row, 18).Value = $clasif
}

# New week (2022-03-08, serial 44628) - rows 304-306.
Set-PrecioRow 304 44628 "Primera" 440 9500  10000 9750  542
Set-PrecioRow 305 44628 "Segunda" 280 7500  8000  7750  431
Set-PrecioRow 306 44628 "Tercera" 200 5500  6000  5750  319

# Previous week (2021-05-19, serial 44335), shifted down - rows 307-310.
Set-PrecioRow 307 44335 "Especial" 600 13500 14000 13750 764
Set-PrecioRow 308 44335 "Primera"  400 11500 12000 11750 653
Set-PrecioRow 309 44335 "Segunda"  300 9500  10000 9750  542
Set-PrecioRow 310 44335 "Tercera"  200 7500  8000  7750  431
